$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.471.62"
$ws.Range("E2").Value = "  -5.99%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.211.01"
$ws.Range("E3").Value = "  -8.91%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.09"
$ws.Range("E5").Value = "  -6.11%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.42"
$ws.Range("E6").Value = "  -13.88%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.205.33"

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -11.56%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -13.84%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  -11.33%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -15.95%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.16"
$ws.Range("E13").Value = "  -17.97%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -13.06%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.719.67"
$ws.Range("E15").Value = "  -9.21%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.438.29"
$ws.Range("E16").Value = "  -6.08%  "

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.209.37"
$ws.Range("E17").Value = "  -8.93%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "533.30"
$ws.Range("E18").Value = "  -13.66%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -6.28%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -15.78%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.98"
$ws.Range("E21").Value = "  -15.60%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.753"
$ws.Range("E22").Value = "  -14.87%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.66"
$ws.Range("E23").Value = "  -14.70%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.07"
$ws.Range("E24").Value = "  -13.60%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.21"
$ws.Range("E25").Value = "  -15.97%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -17.74%  "

# Row 28 - ImmutableX
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -17.97%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  -13.23%  "

# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.93"
$ws.Range("E30").Value = "  -14.48%  "

# Row 31 - Stacks
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.54"
$ws.Range("E31").Value = "  -16.64%  "

# Row 32 - Mantle
$ws.Range("E32").Value = "  -14.26%  "

# Row 33 - Bittensor
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "530.45"
$ws.Range("E33").Value = "  -14.43%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -20.69%  "

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.61"
$ws.Range("E35").Value = "  -18.11%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.17%  "

# Row 37 - OKB
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.01"
$ws.Range("E37").Value = "  -7.11%  "

# Row 38 - now VeChain (was Hedera)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0421"
$ws.Range("E38").Value = "  -13.08%  "

# Row 39 - now Hedera (was VeChain)
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0853"
$ws.Range("E39").Value = "  -15.06%  "

# Row 40 - Cosmos
$ws.Range("E40").Value = "  -16.85%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -14.26%  "

# Row 42 - Maker
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.887.16"
$ws.Range("E42").Value = "  -14.43%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -25.74%  "

# Row 44 - PEPE
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0581"
$ws.Range("E44").Value = "  -21.31%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -17.46%  "

# Row 46 - USDe
$ws.Range("E46").Value = "  -0.07%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  -20.92%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.71"
$ws.Range("E48").Value = "  -20.34%  "

# Row 49 - Fetch.AI
$ws.Range("E49").Value = "  -18.56%  "

# Row 50 - Stellar
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.112"
$ws.Range("E50").Value = "  -14.11%  "

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.74"
$ws.Range("E51").Value = "  -8.18%  "
